# Refresh the cryptocurrency price/volume snapshot (GitHub Actions data pull).
#
# Column D ("Price") in this sheet is always stored as literal text, even when
# it looks like a plain number (e.g. "308.29", "1.000") -- some rows aren't
# even valid numbers at all (thousands-grouped values like "27.173.17"). Excel's
# COM Range.Value setter auto-detects numeric-looking strings and coerces them
# to real numbers, which silently drops formatting such as trailing zeros
# ("1.000" -> 1). To keep every Price cell as text we prefix the assigned
# string with a leading apostrophe, mirroring Excel's own "quote prefix" text
# entry convention; Excel strips the apostrophe and stores the remainder as text.
#
# Rows 39 and 41 swap their Coin/Link (re-ranked upstream) in addition to
# getting refreshed Price/Volume figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''27.170.54'
$ws.Range("E2").Value = '  -2.54%  '

$ws.Range("D3").Value = '''1.712.42'
$ws.Range("E3").Value = '  -3.03%  '

$ws.Range("D4").Value = '''1.001'
$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").Value = '''308.22'
$ws.Range("E5").Value = '  -6.06%  '

$ws.Range("E6").Value = '  +0.04%  '

$ws.Range("D7").Value = '''0.4739'
$ws.Range("E7").Value = '  +6.06%  '

$ws.Range("D8").Value = '''0.3432'
$ws.Range("E8").Value = '  -3.32%  '

$ws.Range("D9").Value = '''42.05'
$ws.Range("E9").Value = '  +0.49%  '

$ws.Range("D10").Value = '''0.07284'
$ws.Range("E10").Value = '  -1.80%  '

$ws.Range("D11").Value = '''1.044'
$ws.Range("E11").Value = '  -5.28%  '

$ws.Range("D12").Value = '''1.001'
$ws.Range("E12").Value = '  +0.04%  '

$ws.Range("D13").Value = '''19.81'
$ws.Range("E13").Value = '  -5.30%  '

$ws.Range("D14").Value = '''5.847'
$ws.Range("E14").Value = '  -3.01%  '

$ws.Range("D15").Value = '''1.713.84'
$ws.Range("E15").Value = '  -2.81%  '

$ws.Range("D16").Value = '''6.839'
$ws.Range("E16").Value = '  -5.54%  '

$ws.Range("D17").Value = '''88.97'
$ws.Range("E17").Value = '  -4.60%  '

$ws.Range("D18").Value = '''0.00001036'
$ws.Range("E18").Value = '  -2.51%  '

$ws.Range("D19").Value = '''0.06352'
$ws.Range("E19").Value = '  -1.23%  '

$ws.Range("D20").Value = '''1.000'
$ws.Range("E20").Value = '  +0.02%  '

$ws.Range("E21").Value = '  -3.67%  '

$ws.Range("D22").Value = '''5.582'
$ws.Range("E22").Value = '  -3.41%  '

$ws.Range("D23").Value = '''27.202.27'
$ws.Range("E23").Value = '  -2.55%  '

$ws.Range("D24").Value = '''10.82'
$ws.Range("E24").Value = '  -3.95%  '

$ws.Range("D25").Value = '''2.095'
$ws.Range("E25").Value = '  -0.85%  '

$ws.Range("D26").Value = '''153.62'
$ws.Range("E26").Value = '  -4.68%  '

$ws.Range("D27").Value = '''19.75'
$ws.Range("E27").Value = '  -3.09%  '

$ws.Range("D28").Value = '''1.901.44'
$ws.Range("E28").Value = '  -3.29%  '

$ws.Range("D29").Value = '''2.089'
$ws.Range("E29").Value = '  -2.96%  '

$ws.Range("D30").Value = '''119.63'
$ws.Range("E30").Value = '  -3.71%  '

$ws.Range("D31").Value = '''1.010'
$ws.Range("E31").Value = '  -8.70%  '

$ws.Range("D32").Value = '''0.09244'
$ws.Range("E32").Value = '  +0.40%  '

$ws.Range("D33").Value = '''3.587'
$ws.Range("E33").Value = '  -2.73%  '

$ws.Range("D34").Value = '''5.289'
$ws.Range("E34").Value = '  -6.66%  '

$ws.Range("D35").Value = '''0.02191'
$ws.Range("E35").Value = '  -3.96%  '

$ws.Range("D36").Value = '''0.05882'
$ws.Range("E36").Value = '  -5.17%  '

$ws.Range("D37").Value = '''11.03'
$ws.Range("E37").Value = '  -7.09%  '

$ws.Range("D38").Value = '''0.1997'
$ws.Range("E38").Value = '  -4.99%  '

$ws.Range("B39").Value = 'WEMIXTOKEN'
$ws.Range("C39").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D39").Value = '''1.412'
$ws.Range("E39").Value = '  +1.33%  '

$ws.Range("D40").Value = '''1.000'
$ws.Range("E40").Value = '  +0.15%  '

$ws.Range("B41").Value = 'InternetComputer(DFINITY)'
$ws.Range("C41").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D41").Value = '''4.733'
$ws.Range("E41").Value = '  -4.76%  '

$ws.Range("D42").Value = '''0.5905'
$ws.Range("E42").Value = '  -6.43%  '

$ws.Range("D43").Value = '''1.111'
$ws.Range("E43").Value = '  -6.14%  '

$ws.Range("D44").Value = '''7.479'
$ws.Range("E44").Value = '  -4.65%  '

$ws.Range("D45").Value = '''12.67'
$ws.Range("E45").Value = '  -4.67%  '

$ws.Range("D46").Value = '''3.560'
$ws.Range("E46").Value = '  -4.93%  '

$ws.Range("D47").Value = '''0.5614'
$ws.Range("E47").Value = '  -4.29%  '

$ws.Range("D48").Value = '''117.96'
$ws.Range("E48").Value = '  -3.59%  '

$ws.Range("D49").Value = '''1.837'
$ws.Range("E49").Value = '  -6.00%  '

$ws.Range("D50").Value = '''0.06628'
$ws.Range("E50").Value = '  -3.63%  '

$ws.Range("D51").Value = '''1.084'
$ws.Range("E51").Value = '  -4.77%  '
